$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "PLN"
$ws.Range("C1").Value = "USD"
$ws.Range("D1").Value = "EUR"
$ws.Range("E1").Value = "CZK"
$ws.Range("F1").Value = "NOK"
$ws.Range("G1").Value = "DKK"

$ws.Range("A1").Value = "currency"
$ws.Range("A2").Value = "amount"

$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = 500
$ws.Range("D2").Value = 1250
$ws.Range("E2").Value = 200
$ws.Range("F2").Value = 600
$ws.Range("G2").Value = 300

$ws.Range("E5").ClearContents()
